# Insert a new data row at row 51 (shifts existing rows 51-109 down to 52-110)
# and populate it with the new "Camote" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(51).Insert()

$ws.Cells.Item(51, 1).Value = 10
$ws.Cells.Item(51, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(51, 3).Value = "La Araucanía"
$ws.Cells.Item(51, 4).Value = 44781
$ws.Cells.Item(51, 5).Value = 9
$ws.Cells.Item(51, 6).Value = 100114002
$ws.Cells.Item(51, 7).Value = "Camote"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 80
$ws.Cells.Item(51, 11).Value = 20000
$ws.Cells.Item(51, 12).Value = 20000
$ws.Cells.Item(51, 13).Value = 20000
$ws.Cells.Item(51, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(51, 15).Value = "Perú"
$ws.Cells.Item(51, 16).Value = 1000
$ws.Cells.Item(51, 17).Value = 20
$ws.Cells.Item(51, 18).Value = "Hortaliza"
